# Task.xlsx edit: populate build/testing dates for rows 8-18, fix a couple
# of "Names" (E column) assignments, and move the active selection to E18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Start/Finish Build dates + "% Build" (column D) for rows 8-16 -------
# (B = Start Build, C = Finish Build, D = % Build; all as date serials /
#  the numeric "1" = 100%, matching the existing numFmtId=164 "mm/dd" style
#  and numFmtId=9 "0%" style already applied to these cells.)

$ws.Range("B8").Value  = 45756
$ws.Range("C8").Value  = 45756
$ws.Range("D8").Value  = 1

$ws.Range("B9").Value  = 45786
$ws.Range("C9").Value  = 45786
$ws.Range("D9").Value  = 1

$ws.Range("B10").Value = 45786
$ws.Range("C10").Value = 45786
$ws.Range("D10").Value = 1

$ws.Range("B11").Value = 45817
$ws.Range("C11").Value = 45817
$ws.Range("D11").Value = 1

$ws.Range("B12").Value = 45847
$ws.Range("C12").Value = 45847
$ws.Range("D12").Value = 1

$ws.Range("B13").Value = 45878
$ws.Range("C13").Value = 45878
$ws.Range("D13").Value = 1

$ws.Range("B14").Value = 45909
$ws.Range("C14").Value = 45939
$ws.Range("D14").Value = 1

$ws.Range("B15").Value = 45970
$ws.Range("C15").Value = 46000
$ws.Range("D15").Value = 1

$ws.Range("B16").Value = 45970
$ws.Range("C16").Value = 46000
$ws.Range("D16").Value = 1

# --- Rows 17-18: dates aren't known yet, so Start/Finish show the literal
#     text "13/09" instead of a real date serial. That needs a right-
#     aligned numFmtId=164 style (distinct from the plain date style
#     already used elsewhere), so set format + alignment explicitly.

$ws.Range("B17").Value = "13/09"
$ws.Range("C17").Value = "13/09"
$ws.Range("D17").Value = 1

$ws.Range("B18").Value = "13/09"
$ws.Range("C18").Value = "13/09"
$ws.Range("D18").Value = 1

$ws.Range("B17:C18").NumberFormat = "mm/dd"
$ws.Range("B17:C18").HorizontalAlignment = -4152   # xlRight

# --- Swap the "Names" assignee for a few rows -----------------------------
$ws.Range("E19").Value = "Hiếu"
$ws.Range("E20").Value = "Hiếu"
$ws.Range("E21").Value = "Hiếu"

$ws.Range("E27").Value = "Hoàng Duy"
$ws.Range("E28").Value = "Hoàng Duy"

# --- Move the active selection -------------------------------------------
$ws.Range("E18").Select() | Out-Null
